$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the M-column totals that were re-measured to a flat 99
$ws.Range("M8").Value  = 99
$ws.Range("M11").Value = 99
$ws.Range("M12").Value = 99
$ws.Range("M17").Value = 99
$ws.Range("M18").Value = 99
$ws.Range("M19").Value = 99
$ws.Range("M20").Value = 99
$ws.Range("M21").Value = 99
$ws.Range("M22").Value = 99
$ws.Range("M23").Value = 99
$ws.Range("M24").Value = 99
$ws.Range("M25").Value = 99

# Update the sheet's view state: the user scrolled down a bit and
# left the selection on S19 instead of Q28
$ws.Range("S19").Select()
